# Monthly EIA update (2017-01-31): add "November" row, update rolling annual
# totals for 2015/2016, and refresh the "as of" month in the title strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_5_07")

# --- 1. Insert a new data row for November at row 50 (old row 50 was the
#        "Rolling 12 Months Ending in ..." banner). ---
# Inserting here shifts the "Rolling 12 months" banner (old 50 -> 51), the
# 2015/2016 total rows (old 51/52 -> 52/53), and the footnote row
# (old 53 -> 54) down by one automatically, including their merged ranges.
$ws.Rows("50").Insert()

# Copy the formatting (styles/borders/number formats) from the October row
# (row 49) down onto the blank inserted row so the new row matches the
# existing monthly-data look exactly.
$ws.Range("A49:F49").Copy()
$ws.Range("A50:F50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Populate the new November row (row 50). ---
$ws.Range("A50").Value = "November"
$ws.Range("B50").Value = 131340346
$ws.Range("C50").Value = 18176753
$ws.Range("D50").Value = 799957
$ws.Range("E50").Value = 83
$ws.Range("F50").Value = 150317139

# --- 3. Refresh the rolling "Year 2015" total (now row 52). ---
$ws.Range("B52").Value = 129711792
$ws.Range("C52").Value = 17972744
$ws.Range("D52").Value = 836010
$ws.Range("E52").Value = 78
$ws.Range("F52").Value = 148520624

# --- 4. Refresh the rolling "Year 2016" total (now row 53). ---
$ws.Range("B53").Value = 131008111
$ws.Range("C53").Value = 18132092
$ws.Range("D53").Value = 815040
$ws.Range("E53").Value = 82
$ws.Range("F53").Value = 149955325

# --- 5. Update the report month references in the title block. ---
$ws.Range("A2").Value = "2008 - November 2016"
$ws.Range("A51").Value = "Rolling 12 Months Ending in November"
